# "location correction for reports"
# Corrects the UNIT SALES TRACKING figures on the studio business report:
#  - adds the week's new Original/Renewal unit sale (row 21, "Units: x.00")
#  - recomputes the Net YTD Total/Sold counts (row 24, "T : n")
#  - recomputes the Net YTD unit totals by category (row 25, "Units: x.00")
#  - recomputes the Net YTD dollar totals by category (row 26, "$x,xxx.xx")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 - this week's unit counts (Pre Original / Original / Extension / Renewal / Total)
$ws.Range("D21").Value = "Units: 2.00"
$ws.Range("H21").Value = "Units: 3.00"
$ws.Range("J21").Value = "Units: 5.00"

# Row 24 - Net YTD Total/Sold (T : n) counts
$ws.Range("D24").Value = "T : 14"
$ws.Range("F24").Value = "T : 11"
$ws.Range("J24").Value = "T : 81"

# Row 25 - Net YTD unit totals by category
$ws.Range("B25").Value = "Units: 166.00"
$ws.Range("D25").Value = "Units: 55.00"
$ws.Range("H25").Value = "Units: 166.00"
$ws.Range("J25").Value = "Units: 547.00"

# Row 26 - Net YTD dollar totals by category
# (leading apostrophe keeps these as literal text, matching the source
#  file where these currency-looking labels are stored as strings)
$ws.Range("B26").Value = "'$20,000.50"
$ws.Range("H26").Value = "'$20,956.90"
$ws.Range("J26").Value = "'$78,187.60"
